# Updates cryptos list prices (D column) and 1h volume % changes (E column)
# for rows 2-51, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.154.68'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '3.320.11'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.19'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.20'
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '3.310.90'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("E10").Value = '  +6.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.632'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.01'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.03'
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("D15").Value = '3.852.77'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.03'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").Value = '3.315.58'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").Value = '64.153.81'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.65'
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.979'
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '449.93'
$ws.Range("E22").Value = '  +5.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.00'
$ws.Range("E23").Value = '  +3.65%  '
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.86'
$ws.Range("E25").Value = '  +3.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.86'
$ws.Range("E26").Value = '  +5.95%  '
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.55'
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.54'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.75'
$ws.Range("E30").Value = '  +3.84%  '
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '62.48'
$ws.Range("E32").Value = '  +7.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.34'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '569.04'
$ws.Range("E34").Value = '  -4.00%  '
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.142'
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.52'
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.08'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("D41").Value = '0.0₃0724'
$ws.Range("E41").Value = '  -3.81%  '
$ws.Range("D42").Value = '3.053.73'
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("E44").Value = '  -3.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.17'
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.59'
$ws.Range("E49").Value = '  +3.24%  '
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("E51").Value = '  -0.27%  '
